$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells that will receive numeric-looking text
# (e.g. "1.000", "5.270") to be formatted as Text first, so Excel
# does not reinterpret/normalize them as numbers and drop meaningful
# trailing zeros / change their textual representation.
$priceRows = @(2,3,5,7,8,10,11,12,13,14,15,17,18,20,21,22,23,24,25,26,27,30,31,32,33,34,35,36,37,38,39,40,41,42,43,45,46,47,50,51)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '30.527.30'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '1.887.31'
$ws.Range("E3").Value = '  +0.75%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '244.29'
$ws.Range("E5").Value = '  -1.46%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.4721'
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '0.2896'
$ws.Range("E8").Value = '  -0.34%  '
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("D10").Value = '22.27'
$ws.Range("E10").Value = '  +1.10%  '
$ws.Range("D11").Value = '0.07762'
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("D12").Value = '1.890.94'
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("D13").Value = '95.92'
$ws.Range("E13").Value = '  -0.13%  '
$ws.Range("D14").Value = '0.7262'
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("D15").Value = '5.191'
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("E16").Value = '  +3.41%  '
$ws.Range("D17").Value = '30.518.55'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = '13.05'
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").Value = '0.000007477'
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").Value = '2.137.13'
$ws.Range("E21").Value = '  +1.00%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '5.270'
$ws.Range("E23").Value = '  +0.55%  '
$ws.Range("D24").Value = '6.335'
$ws.Range("E24").Value = '  +2.68%  '
$ws.Range("D25").Value = '164.50'
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").Value = '9.084'
$ws.Range("E26").Value = '  -1.11%  '
$ws.Range("D27").Value = '18.88'
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("D30").Value = '0.09679'
$ws.Range("E30").Value = '  -2.94%  '
$ws.Range("D31").Value = '1.471'
$ws.Range("E31").Value = '  -2.70%  '
$ws.Range("D32").Value = '4.280'
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("D33").Value = '4.150'
$ws.Range("E33").Value = '  +1.75%  '
$ws.Range("D34").Value = '0.04864'
$ws.Range("E34").Value = '  +1.78%  '
$ws.Range("D35").Value = '1.126'
$ws.Range("E35").Value = '  +0.77%  '
$ws.Range("D36").Value = '0.6933'
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("D37").Value = '2.715'
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").Value = '0.01887'
$ws.Range("E38").Value = '  +1.94%  '
$ws.Range("D39").Value = '2.822'
$ws.Range("E39").Value = '  +2.61%  '
$ws.Range("D40").Value = '74.89'
$ws.Range("E40").Value = '  +2.46%  '
$ws.Range("D41").Value = '6.209'
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("D42").Value = '1.971'
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("D43").Value = '0.4272'
$ws.Range("E43").Value = '  +3.00%  '
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").Value = '0.8274'
$ws.Range("E45").Value = '  -0.72%  '
$ws.Range("D46").Value = '101.31'
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").Value = '9.650'
$ws.Range("E47").Value = '  +3.70%  '
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("D50").Value = '907.18'
$ws.Range("E50").Value = '  -0.51%  '
$ws.Range("D51").Value = '0.05748'
$ws.Range("E51").Value = '  +1.62%  '
